$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Televisor PANASONIC LED 32'''' HD Smart TV T...'
$ws.Range("B1").Value = 'S/. 899.00'
$ws.Range("C1").Value = 'S/0'
$ws.Range("D1").Value = 'S/0'
$ws.Range("A2").Value = 'Televisor Xiaomi Mi LED TV 4A 32'' Smart HD...'
$ws.Range("B2").Value = 'S/. 989.00'
$ws.Range("A3").Value = 'Televisor HYUNDAI LED 58" UHD 4K Smart Tv ...'
$ws.Range("B3").Value = 'S/. 1,899.00'
$ws.Range("C3").Value = 'S/0'
$ws.Range("D3").Value = 'S/0'
$ws.Range("A4").Value = 'Televisor LG HD ThinQ AI 32" 32LM637B'
$ws.Range("B4").Value = 'S/0'
$ws.Range("A5").Value = 'Televisor LG NanoCell 4K Procesador Inteli...'
$ws.Range("B5").Value = 'S/0'
$ws.Range("C5").Value = 'S/0'
$ws.Range("D5").Value = 'S/0'
$ws.Range("A6").Value = 'Televisor LG NanoCell 4K ThinQ AI 70" 70NA...'
$ws.Range("B6").Value = 'S/. 4,299.00'
$ws.Range("A7").Value = 'TV Smart LG 4K 75" NanoCell, Thinq Ai, Ult...'
$ws.Range("B7").Value = 'S/. 4,999.00'
$ws.Range("A8").Value = 'TELEVISOR LG 75" UHD 4K MOD: 75UP7750PSB'
$ws.Range("B8").Value = 'S/. 4,399.00'
$ws.Range("A9").Value = 'Xiaomi TV Smart 43" 4K UHD Modelo: L43M6'
$ws.Range("B9").Value = 'S/. 1,499.00'
$ws.Range("A10").Value = 'Televisor Hisense LED UHD 58" 58A6GSV'
$ws.Range("C10").Value = 'S/0'
$ws.Range("D10").Value = 'S/0'
$ws.Range("A11").Value = 'Televisor Hisense LED HD 32" 32A4GSV'
$ws.Range("B11").Value = 'S/. 749.00'
$ws.Range("C11").Value = 'S/0'
$ws.Range("D11").Value = 'S/0'
$ws.Range("A12").Value = 'Televisor Hisense LED UHD 50" 50A6GSV'
$ws.Range("B12").Value = 'S/. 1,499.00'
$ws.Range("A13").Value = 'COMBO Televisor LG Smart TV UHD 55" + Máqu...'
$ws.Range("B13").Value = 'S/. 2,099.00'
$ws.Range("C13").Value = 'S/0'
$ws.Range("D13").Value = 'S/0'
$ws.Range("A14").Value = 'TV Smart Xiaomi Mi TV P1 32" LED, HD, sist...'
$ws.Range("B14").Value = 'S/0'
$ws.Range("A15").Value = 'TV Smart Xiaomi Mi TV Q1 4K 75" QLED, Ultr...'
$ws.Range("B15").Value = 'S/0'
$ws.Range("A16").Value = 'TV Smart Xiaomi Mi TV P1 4K 50" LED, Ultra...'
$ws.Range("B16").Value = 'S/0'
$ws.Range("A17").Value = 'TV Smart Xiaomi Mi TV P1 4K 55" LED, Ultra...'
$ws.Range("A18").Value = 'Soporte de Pared Fijo Inclinable para TV 3...'
$ws.Range("B18").Value = 'S/0'
$ws.Range("A19").Value = 'Televisor LED SMART HD 32" 32S5195'
$ws.Range("B19").Value = 'S/. 949.00'
$ws.Range("A20").Value = 'Televisor AOC LED SMART FHD 43" 43S5195'
$ws.Range("B20").Value = 'S/. 1,399.00'
$ws.Range("A21").Value = 'TELEVISOR SAMSUNG 65" MOD: QN65Q70AAGXPE'
$ws.Range("B21").Value = 'S/. 5,179.00'
$ws.Range("A22").Value = 'TELEVISOR SAMSUNG 65" MOD: UN65AU8000GXPE'
$ws.Range("B22").Value = 'S/. 3,579.00'
$ws.Range("A23").Value = 'TELEVISOR SAMSUNG 55" MOD: UN55AU8000GXPE'
$ws.Range("B23").Value = 'S/. 2,539.00'
$ws.Range("A24").Value = 'TELEVISOR SAMSUNG 58" MOD: UN58AU7000GXPE'
$ws.Range("B24").Value = 'S/. 2,539.00'
